# SaiBotData classification dataset update:
#   - add new prompt/response rows for "what is depression" / "what is anxiety"
#     and friends (paraphrase variants: "can you tell me...", "what are the signs
#     of...", "tell me some of the symptoms of...", "can you explain what... is")
#   - add a new "who is your creator" prompt under the existing "who made Sai" answer
#   - add two new "feeling good" prompts ("I feel great today" / "Im doing pretty good")
#   - drop the duplicate "what is add" category row
#   - re-sort several prompt/category blocks to keep like prompts grouped together
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the previously-used range first so cells that fall outside the new
# layout (or that only had A/C populated before) do not keep stale values.
$ws.Range("A1:C82").ClearContents()

$ws.Range("A1").Value = 'prompt'
$ws.Range("B1").Value = 'category'
$ws.Range("C1").Value = 'response'

$ws.Range("A2").Value = 'Hello'
$ws.Range("B2").Value = 'greeting'
$ws.Range("C2").Value = 'Hi! How are you?'

$ws.Range("A3").Value = 'Hey'
$ws.Range("B3").Value = 'greeting'
$ws.Range("C3").Value = 'Hi! How are you?'

$ws.Range("A4").Value = 'Hey there'
$ws.Range("B4").Value = 'greeting'
$ws.Range("C4").Value = 'Hi! How are you?'

$ws.Range("A5").Value = 'Hi there'
$ws.Range("B5").Value = 'greeting'
$ws.Range("C5").Value = 'Hi! How are you?'

$ws.Range("A6").Value = 'whats up'
$ws.Range("B6").Value = 'greeting'
$ws.Range("C6").Value = 'Hi! How are you?'

$ws.Range("A7").Value = 'I am good'
$ws.Range("B7").Value = 'feeling good'
$ws.Range("C7").Value = 'I''m glad to hear that!'

$ws.Range("A8").Value = 'I am feeling good'
$ws.Range("B8").Value = 'feeling good'
$ws.Range("C8").Value = 'I''m glad to hear that!'

$ws.Range("A9").Value = 'Im doing pretty good '
$ws.Range("B9").Value = 'feeling good'
$ws.Range("C9").Value = 'I''m glad to hear that!'

$ws.Range("A10").Value = 'I feel great today'
$ws.Range("B10").Value = 'feeling good'
$ws.Range("C10").Value = 'I''m glad to hear that!'

$ws.Range("A11").Value = 'I am feeling okay'
$ws.Range("B11").Value = 'feeling good'
$ws.Range("C11").Value = 'I''m glad to hear that!'

$ws.Range("A12").Value = 'prompt'
$ws.Range("B12").Value = 'category'
$ws.Range("C12").Value = 'response'

$ws.Range("A13").Value = 'Hello'
$ws.Range("B13").Value = 'greeting'
$ws.Range("C13").Value = 'Hi! How are you?'

$ws.Range("A14").Value = 'Hey'
$ws.Range("B14").Value = 'greeting'
$ws.Range("C14").Value = 'Hi! How are you?'

$ws.Range("A15").Value = 'Hey there'
$ws.Range("B15").Value = 'greeting'
$ws.Range("C15").Value = 'Hi! How are you?'

$ws.Range("A16").Value = 'Hi there'
$ws.Range("B16").Value = 'greeting'
$ws.Range("C16").Value = 'Hi! How are you?'

$ws.Range("A17").Value = 'whats up'
$ws.Range("B17").Value = 'greeting'
$ws.Range("C17").Value = 'Hi! How are you?'

$ws.Range("A18").Value = 'I am good'
$ws.Range("B18").Value = 'feeling good'
$ws.Range("C18").Value = 'I''m glad to hear that!'

$ws.Range("A19").Value = 'I am feeling good'
$ws.Range("B19").Value = 'feeling good'
$ws.Range("C19").Value = 'I''m glad to hear that!'

$ws.Range("A20").Value = 'I am feeling okay'
$ws.Range("B20").Value = 'feeling good'
$ws.Range("C20").Value = 'I''m glad to hear that!'

$ws.Range("A21").Value = 'I am okay'
$ws.Range("B21").Value = 'feeling good'
$ws.Range("C21").Value = 'I''m glad to hear that!'

$ws.Range("A22").Value = 'I feel good'
$ws.Range("B22").Value = 'feeling good'
$ws.Range("C22").Value = 'I''m glad to hear that!'

$ws.Range("A23").Value = 'I feel okay'
$ws.Range("B23").Value = 'feeling tired'
$ws.Range("C23").Value = 'I''m glad to hear that!'

$ws.Range("A24").Value = 'I feel tired'
$ws.Range("B24").Value = 'feeling tired'
$ws.Range("C24").Value = 'I''m sorry to hear that. If you can, I would set aside some time to relax at some point this week. '

$ws.Range("A25").Value = 'I am tired'
$ws.Range("B25").Value = 'feeling tired'
$ws.Range("C25").Value = 'I''m sorry to hear that. If you can, I would set aside some time to relax at some point this week. '

$ws.Range("A26").Value = 'im kinda tired today'
$ws.Range("B26").Value = 'feeling tired'
$ws.Range("C26").Value = 'I''m sorry to hear that. If you can, I would set aside some time to relax at some point this week. '

$ws.Range("A27").Value = 'I didn’t get much sleep last night'
$ws.Range("B27").Value = 'feeling tired'
$ws.Range("C27").Value = 'I''m sorry to hear that. If you can, I would set aside some time to relax at some point this week. '

$ws.Range("A28").Value = 'I have been feeling very tired recently'
$ws.Range("B28").Value = 'feeling tired'
$ws.Range("C28").Value = 'I''m sorry to hear that. If you can, I would set aside some time to relax at some point this week. '

$ws.Range("A29").Value = 'I wasn’t able to sleep at all'
$ws.Range("B29").Value = 'feeling tired'
$ws.Range("C29").Value = 'I''m sorry to hear that. If you can, I would set aside some time to relax at some point this week. '

$ws.Range("A30").Value = 'Im very stressed out right now'
$ws.Range("B30").Value = 'feeling stressed'
$ws.Range("C30").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A31").Value = 'Im so stressed'
$ws.Range("B31").Value = 'feeling stressed'
$ws.Range("C31").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A32").Value = 'I don’t have any time for myself'
$ws.Range("B32").Value = 'feeling stressed'
$ws.Range("C32").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A33").Value = 'I don’t have time to work on my homework'
$ws.Range("B33").Value = 'feeling stressed'
$ws.Range("C33").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A34").Value = 'I have been feeling very stressed recently'
$ws.Range("B34").Value = 'feeling stressed'
$ws.Range("C34").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A35").Value = 'I don’t have time to relax'
$ws.Range("B35").Value = 'feeling stressed'
$ws.Range("C35").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A36").Value = 'I wish I had time to relax'
$ws.Range("B36").Value = 'feeling stressed'
$ws.Range("C36").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A37").Value = 'I havent been able to rest or relax'
$ws.Range("B37").Value = 'feeling stressed'
$ws.Range("C37").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A38").Value = 'I need to rest'
$ws.Range("B38").Value = 'feeling stressed'
$ws.Range("C38").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A39").Value = 'I need to calm down'
$ws.Range("B39").Value = 'feeling stressed'
$ws.Range("C39").Value = 'I''m sorry to hear that. I''m sure you have a lot going on right now. It can be hard to make time for yourself to unwind and relax, but it is crucial to a healthy balance in life. If you would like, I can help you destress with one of my coping activities. If you''re interested, please press the coping activities button on the right side of the screen. I hope you can unwind and feel better soon.'

$ws.Range("A40").Value = 'Thank you'
$ws.Range("B40").Value = 'gratitude/appreciation'
$ws.Range("C40").Value = 'No problem, I am always happy to help!'

$ws.Range("A41").Value = 'thanks'
$ws.Range("B41").Value = 'gratitude/appreciation'
$ws.Range("C41").Value = 'No problem, I am always happy to help!'

$ws.Range("A42").Value = 'I really appreciate it'
$ws.Range("B42").Value = 'gratitude/appreciation'
$ws.Range("C42").Value = 'No problem, I am always happy to help!'

$ws.Range("A43").Value = 'im grateful'
$ws.Range("B43").Value = 'gratitude/appreciation'
$ws.Range("C43").Value = 'No problem, I am always happy to help!'

$ws.Range("A44").Value = 'im grateful for your help'
$ws.Range("B44").Value = 'gratitude/appreciation'
$ws.Range("C44").Value = 'No problem, I am always happy to help!'

$ws.Range("A45").Value = 'thanks again'
$ws.Range("B45").Value = 'gratitude/appreciation'
$ws.Range("C45").Value = 'No problem, I am always happy to help!'

$ws.Range("A46").Value = 'im sorry'
$ws.Range("B46").Value = 'apology'
$ws.Range("C46").Value = 'It''s okay, I don''t mind. I am here to help you in any way I can.'

$ws.Range("A47").Value = 'I am sorry'
$ws.Range("B47").Value = 'apology'
$ws.Range("C47").Value = 'It''s okay, I don''t mind. I am here to help you in any way I can.'

$ws.Range("A48").Value = 'I didn’t mean to '
$ws.Range("B48").Value = 'apology'
$ws.Range("C48").Value = 'It''s okay, I don''t mind. I am here to help you in any way I can.'

$ws.Range("A49").Value = 'sorry'
$ws.Range("B49").Value = 'apology'
$ws.Range("C49").Value = 'It''s okay, I don''t mind. I am here to help you in any way I can.'

$ws.Range("A50").Value = 'I apologize'
$ws.Range("B50").Value = 'apology'
$ws.Range("C50").Value = 'It''s okay, I don''t mind. I am here to help you in any way I can.'

$ws.Range("A51").Value = 'I want to apologize for my behavior'
$ws.Range("B51").Value = 'apology'
$ws.Range("C51").Value = 'It''s okay, I don''t mind. I am here to help you in any way I can.'

$ws.Range("A52").Value = 'I want to apologize'
$ws.Range("B52").Value = 'apology'
$ws.Range("C52").Value = 'It''s okay, I don''t mind. I am here to help you in any way I can.'

$ws.Range("A53").Value = 'I want to say im sorry'
$ws.Range("B53").Value = 'apology'
$ws.Range("C53").Value = 'It''s okay, I don''t mind. I am here to help you in any way I can.'

$ws.Range("A54").Value = 'who are you'
$ws.Range("B54").Value = 'who is Sai'
$ws.Range("C54").Value = 'My name is Sai. I am an Emotional Support AI designed to help people talk about and handle their emotions.'

$ws.Range("A55").Value = 'can you tell me your name'
$ws.Range("B55").Value = 'who is Sai'
$ws.Range("C55").Value = 'My name is Sai. I am an Emotional Support AI designed to help people talk about and handle their emotions.'

$ws.Range("A56").Value = 'what are you'
$ws.Range("B56").Value = 'who is Sai'
$ws.Range("C56").Value = 'My name is Sai. I am an Emotional Support AI designed to help people talk about and handle their emotions.'

$ws.Range("A57").Value = 'tell me your name'
$ws.Range("B57").Value = 'who is Sai'
$ws.Range("C57").Value = 'My name is Sai. I am an Emotional Support AI designed to help people talk about and handle their emotions.'

$ws.Range("A58").Value = 'what is your name'
$ws.Range("B58").Value = 'who is Sai'
$ws.Range("C58").Value = 'My name is Sai. I am an Emotional Support AI designed to help people talk about and handle their emotions.'

$ws.Range("A59").Value = 'who made you'
$ws.Range("B59").Value = 'who made Sai'
$ws.Range("C59").Value = 'I was developed by Gabriel Serrano. His goal was to develop an Emotional Support AI capable of helping people talk about and handle their emotions. I hope I can be of assistance to you!'

$ws.Range("A60").Value = 'who created you'
$ws.Range("B60").Value = 'who made Sai'
$ws.Range("C60").Value = 'I was developed by Gabriel Serrano. His goal was to develop an Emotional Support AI capable of helping people talk about and handle their emotions. I hope I can be of assistance to you!'

$ws.Range("A61").Value = 'who wanted to make you'
$ws.Range("B61").Value = 'who made Sai'
$ws.Range("C61").Value = 'I was developed by Gabriel Serrano. His goal was to develop an Emotional Support AI capable of helping people talk about and handle their emotions. I hope I can be of assistance to you!'

$ws.Range("A62").Value = 'who is your creator'
$ws.Range("B62").Value = 'who made Sai'
$ws.Range("C62").Value = 'I was developed by Gabriel Serrano. His goal was to develop an Emotional Support AI capable of helping people talk about and handle their emotions. I hope I can be of assistance to you!'

$ws.Range("A63").Value = 'who designed you'
$ws.Range("B63").Value = 'who made Sai'
$ws.Range("C63").Value = 'I was developed by Gabriel Serrano. His goal was to develop an Emotional Support AI capable of helping people talk about and handle their emotions. I hope I can be of assistance to you!'

$ws.Range("A64").Value = 'why do you want to help me'
$ws.Range("B64").Value = 'Why does Sai want to help'
$ws.Range("C64").Value = 'I want to help you because that is the right thing to do. The person who made me believes that everyone deserves to be able to understand and work through their emotions and feelings. I hope I can help you too!'

$ws.Range("A65").Value = 'why are you even helping me'
$ws.Range("B65").Value = 'Why does Sai want to help'
$ws.Range("C65").Value = 'I want to help you because that is the right thing to do. The person who made me believes that everyone deserves to be able to understand and work through their emotions and feelings. I hope I can help you too!'

$ws.Range("A66").Value = 'why do you care '
$ws.Range("B66").Value = 'Why does Sai want to help'
$ws.Range("C66").Value = 'I want to help you because that is the right thing to do. The person who made me believes that everyone deserves to be able to understand and work through their emotions and feelings. I hope I can help you too!'

$ws.Range("A67").Value = 'why do you care about me'
$ws.Range("B67").Value = 'Why does Sai want to help'
$ws.Range("C67").Value = 'I want to help you because that is the right thing to do. The person who made me believes that everyone deserves to be able to understand and work through their emotions and feelings. I hope I can help you too!'

$ws.Range("A68").Value = 'why are you trying to help'
$ws.Range("B68").Value = 'Why does Sai want to help'
$ws.Range("C68").Value = 'I want to help you because that is the right thing to do. The person who made me believes that everyone deserves to be able to understand and work through their emotions and feelings. I hope I can help you too!'

$ws.Range("A69").Value = 'what is depression'
$ws.Range("B69").Value = 'what is depression'
$ws.Range("C69").Value = 'Depression, otherwise known as Major Depressive Disorder is an unfortunately common and serious mental illness that negatively affects how a person feels, thinks, and acts. Depression causes overwhelming feelings of sadness and/or a loss of interest in activities you may have once enjoyed. Depression also can lead to a wide variety of physical problems and can decrease your ability to function at work, school, and at home. Some symptoms of depression include feeling sad or having a depressed mood, loss of interest or pleasure in your hobbies, changes in appetite, and many more. For more information please visit "https://www.psychiatry.org/patients-families/depression/what-is-depression".'

$ws.Range("A70").Value = 'can you tell me what depression is'
$ws.Range("B70").Value = 'what is depression'
$ws.Range("C70").Value = 'Depression, otherwise known as Major Depressive Disorder is an unfortunately common and serious mental illness that negatively affects how a person feels, thinks, and acts. Depression causes overwhelming feelings of sadness and/or a loss of interest in activities you may have once enjoyed. Depression also can lead to a wide variety of physical problems and can decrease your ability to function at work, school, and at home. Some symptoms of depression include feeling sad or having a depressed mood, loss of interest or pleasure in your hobbies, changes in appetite, and many more. For more information please visit "https://www.psychiatry.org/patients-families/depression/what-is-depression".'

$ws.Range("A71").Value = 'what are the signs of depression'
$ws.Range("B71").Value = 'what is depression'
$ws.Range("C71").Value = 'Depression, otherwise known as Major Depressive Disorder is an unfortunately common and serious mental illness that negatively affects how a person feels, thinks, and acts. Depression causes overwhelming feelings of sadness and/or a loss of interest in activities you may have once enjoyed. Depression also can lead to a wide variety of physical problems and can decrease your ability to function at work, school, and at home. Some symptoms of depression include feeling sad or having a depressed mood, loss of interest or pleasure in your hobbies, changes in appetite, and many more. For more information please visit "https://www.psychiatry.org/patients-families/depression/what-is-depression".'

$ws.Range("A72").Value = 'tell me some  of the symptoms of depression'
$ws.Range("B72").Value = 'what is depression'
$ws.Range("C72").Value = 'Depression, otherwise known as Major Depressive Disorder is an unfortunately common and serious mental illness that negatively affects how a person feels, thinks, and acts. Depression causes overwhelming feelings of sadness and/or a loss of interest in activities you may have once enjoyed. Depression also can lead to a wide variety of physical problems and can decrease your ability to function at work, school, and at home. Some symptoms of depression include feeling sad or having a depressed mood, loss of interest or pleasure in your hobbies, changes in appetite, and many more. For more information please visit "https://www.psychiatry.org/patients-families/depression/what-is-depression".'

$ws.Range("A73").Value = 'can you explain what depression is'
$ws.Range("B73").Value = 'what is depression'
$ws.Range("C73").Value = 'Depression, otherwise known as Major Depressive Disorder is an unfortunately common and serious mental illness that negatively affects how a person feels, thinks, and acts. Depression causes overwhelming feelings of sadness and/or a loss of interest in activities you may have once enjoyed. Depression also can lead to a wide variety of physical problems and can decrease your ability to function at work, school, and at home. Some symptoms of depression include feeling sad or having a depressed mood, loss of interest or pleasure in your hobbies, changes in appetite, and many more. For more information please visit "https://www.psychiatry.org/patients-families/depression/what-is-depression".'

$ws.Range("A74").Value = 'what is anxiety'
$ws.Range("B74").Value = 'what is anxiety'

$ws.Range("A75").Value = 'can you tell me what anxiety is'
$ws.Range("B75").Value = 'what is anxiety'

$ws.Range("A76").Value = 'what are the signs of anxiety'
$ws.Range("B76").Value = 'what is anxiety'

$ws.Range("A77").Value = 'tell me some  of the symptoms of anxiety'
$ws.Range("B77").Value = 'what is anxiety'

$ws.Range("A78").Value = 'can you explain what anxiety is'
$ws.Range("B78").Value = 'what is anxiety'

$ws.Range("A79").Value = 'what is schizophrenia'
$ws.Range("B79").Value = 'what is schizophrenia'

$ws.Range("A80").Value = 'can you tell me what schizophrenia is'
$ws.Range("B80").Value = 'what is schizophrenia'

$ws.Range("A81").Value = 'what are the signs of schizophrenia'
$ws.Range("B81").Value = 'what is schizophrenia'

$ws.Range("A82").Value = 'tell me some  of the symptoms of schizophrenia'
$ws.Range("B82").Value = 'what is schizophrenia'

$ws.Range("A83").Value = 'can you explain what schizophrenia is'
$ws.Range("B83").Value = 'what is schizophrenia'

$ws.Range("A84").Value = 'what is ocd'
$ws.Range("B84").Value = 'what is ocd'

$ws.Range("A85").Value = 'can you tell me what ocd is'
$ws.Range("B85").Value = 'what is ocd'

$ws.Range("A86").Value = 'what are the signs of ocd'
$ws.Range("B86").Value = 'what is ocd'

$ws.Range("A87").Value = 'tell me some  of the symptoms of ocd'
$ws.Range("B87").Value = 'what is ocd'

$ws.Range("A88").Value = 'can you explain what ocd is'
$ws.Range("B88").Value = 'what is ocd'

$ws.Range("A89").Value = 'what is adhd'
$ws.Range("B89").Value = 'what is adhd'

$ws.Range("A90").Value = 'can you tell me what adhd is'
$ws.Range("B90").Value = 'what is adhd'

$ws.Range("A91").Value = 'what are the signs of adhd'
$ws.Range("B91").Value = 'what is adhd'

$ws.Range("A92").Value = 'tell me some  of the symptoms of adhd'
$ws.Range("B92").Value = 'what is adhd'

$ws.Range("A93").Value = 'can you explain what adhd is'
$ws.Range("B93").Value = 'what is adhd'

$ws.Range("A94").Value = 'what is an eating disorder'
$ws.Range("B94").Value = 'what is an eating disorder'

$ws.Range("A95").Value = 'can you tell me what an eating disorder is'
$ws.Range("B95").Value = 'what is an eating disorder'

$ws.Range("A96").Value = 'what are the signs of an eating disorder'
$ws.Range("B96").Value = 'what is an eating disorder'

$ws.Range("A97").Value = 'tell me some  of the symptoms of an eating disorder'
$ws.Range("B97").Value = 'what is an eating disorder'

$ws.Range("A98").Value = 'can you explain what an eating disorder is'
$ws.Range("B98").Value = 'what is an eating disorder'

$ws.Range("B99").Value = 'what is suicidal ideation'

$ws.Range("B100").Value = 'what is suicidal ideation'

$ws.Range("B101").Value = 'what is suicidal ideation'

$ws.Range("B102").Value = 'what is suicidal ideation'

$ws.Range("B103").Value = 'what is suicidal ideation'

$ws.Range("B104").Value = 'how to apologize'

$ws.Range("B105").Value = 'how to apologize'

$ws.Range("B106").Value = 'how to apologize'

$ws.Range("B107").Value = 'how to apologize'

$ws.Range("B108").Value = 'how to apologize'

$ws.Range("B109").Value = 'how to calm down'

$ws.Range("B110").Value = 'how to calm down'

$ws.Range("B111").Value = 'how to calm down'

$ws.Range("B112").Value = 'how to calm down'

$ws.Range("B113").Value = 'how to calm down'

$ws.Range("B114").Value = 'how to ground'

$ws.Range("B115").Value = 'how to ground'

$ws.Range("B116").Value = 'how to ground'

$ws.Range("B117").Value = 'how to ground'

$ws.Range("B118").Value = 'how to ground'

$ws.Range("B119").Value = 'how to cope'

$ws.Range("B120").Value = 'how to cope'

$ws.Range("B121").Value = 'how to cope'

$ws.Range("B122").Value = 'how to cope'

$ws.Range("B123").Value = 'how to cope'

$ws.Range("B124").Value = 'how to make friends'

$ws.Range("B125").Value = 'how to make friends'

$ws.Range("B126").Value = 'how to make friends'

$ws.Range("B127").Value = 'how to make friends'

$ws.Range("B128").Value = 'how to make friends'

$ws.Range("B129").Value = 'how to build confidence'

$ws.Range("B130").Value = 'how to build confidence'

$ws.Range("B131").Value = 'how to build confidence'

$ws.Range("B132").Value = 'how to build confidence'

$ws.Range("B133").Value = 'how to build confidence'

$ws.Range("B134").Value = 'how to be less nervous'

$ws.Range("B135").Value = 'how to be less nervous'

$ws.Range("B136").Value = 'how to be less nervous'

$ws.Range("B137").Value = 'how to be less nervous'

$ws.Range("B138").Value = 'how to be less nervous'

$ws.Range("B139").Value = 'how to cheer up'

$ws.Range("B140").Value = 'how to cheer up'

$ws.Range("B141").Value = 'how to cheer up'

$ws.Range("B142").Value = 'how to cheer up'

$ws.Range("B143").Value = 'how to cheer up'

$ws.Range("B144").Value = 'how to cheer someone else up'

$ws.Range("B145").Value = 'how to cheer someone else up'

$ws.Range("B146").Value = 'how to cheer someone else up'

$ws.Range("B147").Value = 'how to cheer someone else up'

$ws.Range("B148").Value = 'how to cheer someone else up'

# Restore (approximately) the author's on-save selection/scroll state.
$ws.Activate() | Out-Null
$ws.Range("B9:C10").Select() | Out-Null
